# Apply the "ajout de données de sortie" edit:
#  - insert a new "Closing1d" column after Price (old G shifts -> H, etc.)
#  - add BDH(...) lookup formulas in the new column for Roll/Outright rows
#  - add B-column Level formula (=Price/Closing1d) for Outright rows
#  - format the Level column (B) as 0.000
#  - clear the stray empty Price cells on the "Screen" rows
#  - keep the Date column formatted as YYYY-MM-DD

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at G ("Closing1d"), shifting Volume..OpenInt right by one.
$ws.Columns.Item(7).Insert()
$ws.Range("G1").Value = "Closing1d"

# 2. Rows that are a "Roll" or an "Outright" get a BDH closing-price lookup in
#    the new column (UndTkr is now column J, Date is now column O).
$bdhRows = @(2,5,8,11,14,17,20,23,26,31,32,33,34,35,36,37,38,39,40,41,42,43)
foreach ($r in $bdhRows) {
    $ws.Cells.Item($r, 7).Formula = '=BDH(J' + $r + '&" Index", "PX_CLOSE_1D",O' + $r + ',O' + $r + ')'
}

# 3. "Outright" rows compute Level as Price / Closing1d instead of a static value.
$outrightRows = @(31,32,33,34,35,36,37,38,39,40,41,42,43)
foreach ($r in $outrightRows) {
    $ws.Cells.Item($r, 2).Formula = '=F' + $r + '/G' + $r
}

# 4. The "Screen" rows never had a real Price - drop the stray empty cell.
$screenRows = @(27,28,29,30)
foreach ($r in $screenRows) {
    $ws.Range("F" + $r).ClearContents()
}

# 5. Roll rows never had an OpenInt value (just a stray empty cell) - drop it.
$emptyOpenIntRows = @(2,5,8,11,14,17,20,23)
foreach ($r in $emptyOpenIntRows) {
    $ws.Range("S" + $r).ClearContents()
}

# 6. Level column (B) is now a ratio - show it with 3 decimals.
$ws.Range("B2:B43").NumberFormat = "0.000"

# 7. Date column (now O, was N) keeps its YYYY-MM-DD display format.
$ws.Range("O2:O43").NumberFormat = "YYYY-MM-DD"
